# DOC: Add future work. XLS: Add PE
#
# The workbook starts with a single sheet ("Sheet1") holding the area
# estimate template. The edit:
#   1. Renames that sheet to "Manager" and tweaks its "H31" line item from
#      a formula to a flat (hardcoded) value, which ripples through the
#      totals below it.
#   2. Duplicates "Manager" to a new "PE" tab (placed right after
#      Manager) representing another block, and adjusts a handful of its
#      numbers (its own H31->H38 total is overridden by hand, and a
#      couple of the "Available Area" inputs are zeroed out).
#   3. Adds a brand-new small "Sheet3" summarizing Manager + PE together,
#      and leaves that as the active/selected tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename the original (only) sheet to "Manager".
# ---------------------------------------------------------------------
$mgr = $wb.Worksheets.Item(1)
$mgr.Name = "Manager"

# ---------------------------------------------------------------------
# 2. Duplicate it right after itself and rename the copy to "PE".
# ---------------------------------------------------------------------
$mgr.Copy([System.Reflection.Missing]::Value, $mgr)
$pe = $wb.Worksheets.Item(2)
$pe.Name = "PE"

# ---------------------------------------------------------------------
# 3. Add the new summary sheet after PE.
# ---------------------------------------------------------------------
$sheet3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $pe)
$sheet3.Name = "Sheet3"

# ---------------------------------------------------------------------
# Manager edits: the "Stack Up" line (H31) becomes a flat value instead
# of F31*G31; everything below it (H38/H40/H42, and the Q19/Q22 ratios
# up top that depend on H42) recalculates automatically.
# ---------------------------------------------------------------------
$mgr.Range("H31").Value = 350000

# ---------------------------------------------------------------------
# PE edits: its own H31 line is also hand-entered, its "Available Area"
# sub-items N34/N35 are zeroed, its H38 sub-total is overridden by hand
# (no longer a SUM formula), and F39 drops to 0. All downstream totals
# (H39/H40/H41/H42, N36/N38/N41/N42/F41, Q19/Q22) recalc automatically.
# ---------------------------------------------------------------------
$pe.Range("H31").Value = 350000
$pe.Range("N34").Value = 0
$pe.Range("N35").Value = 0
$pe.Range("F39").Value = 0
$pe.Range("H38").Value = 2250000

# ---------------------------------------------------------------------
# Sheet3: small roll-up comparing Manager and PE.
# ---------------------------------------------------------------------
$sheet3.Range("H4").Value = 1
$sheet3.Range("E5").Formula = "=Manager!H40"
$sheet3.Range("F5").Formula = "=PE!H40"
$sheet3.Range("G5").Formula = "=SUM(E5:F5)"
$sheet3.Range("H5").Formula = "=G5*H`$4"
$sheet3.Range("E6").Formula = "=E5/Manager!O19"
$sheet3.Range("F6").Formula = "=F5/PE!O19"
$sheet3.Range("G6").Formula = "=SUM(E6:F6)"
$sheet3.Range("H6").Formula = "=G6*H`$4"

# ---------------------------------------------------------------------
# View state: restore each sheet's scroll position / selection, and
# leave Sheet3 as the active tab (matching the authored workbook).
# ---------------------------------------------------------------------
$mgr.Activate()
$excel.ActiveWindow.ScrollRow = 17
$excel.ActiveWindow.ScrollColumn = 1
$mgr.Range("N19").Select()

$pe.Activate()
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 1
$pe.Range("H45").Select()

$sheet3.Activate()
$sheet3.Range("J4").Select()
